$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..14 (Student ID, Subject, Log Date, Log Time, Type, User)
$rows = @(
    @("211177","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("191088","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("211169","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("211741","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("201529","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("211245","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("200359","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("201197","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("201218","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("200852","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("202004","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("200405","general surgery","25/10/2025","10:30:00","Excuse","System"),
    @("211174","general surgery","25/10/2025","10:30:00","Excuse","System")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $rowVals = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
}
